$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new calibration row at row 4 for "2,4,5-trichlorophenol"
# (the sheet is sorted by MW, and this compound's MW=197.4 sorts ahead of
# tetradecanoic acid's 228.37, so it belongs right after row 3).
# Inserting the row copies the formatting down from the row that used to be
# row 4, which is already a very close match for the rest of the table.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# --- Populate the new row with its calibration data ------------------------
$ws.Range("A4").Value = "2,4,5-trichlorophenol"
$ws.Range("B4").Value = 197.4
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = "n.a."
$ws.Range("H4").Value = "n.a."
$ws.Range("I4").Value = 155710
$ws.Range("J4").Value = 343277
$ws.Range("K4").Value = 805095
$ws.Range("L4").Value = 2302730
$ws.Range("M4").Value = "n.a."
$ws.Range("N4").Value = "n.a."

# --- Formatting tweaks so the new row matches the rest of the table --------
# Name cell: left/top aligned, non-themed font (matches the other "pasted in"
# compound name cell style used throughout the workbook).
$ws.Range("A4").HorizontalAlignment = -4131   # xlLeft
$ws.Range("A4").VerticalAlignment = -4160     # xlTop
$ws.Range("A4").Font.Name = "Calibri"

# MW cell: numeric, right aligned, no vertical centering.
$ws.Range("B4").NumberFormat = "0.00"
$ws.Range("B4").HorizontalAlignment = -4152   # xlRight
$ws.Range("B4").VerticalAlignment = -4107     # xlBottom

# PPM6 cell keeps the thin right border used by the rest of the PPM6 column.
$ws.Range("H2").Copy()
$ws.Range("H4").PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = 0

# Area columns: right aligned, no vertical centering (matches Area1-4 style
# used elsewhere in the table).
$ws.Range("I4:L4").HorizontalAlignment = -4152
$ws.Range("I4:L4").VerticalAlignment = -4107

# --- Update the hidden filter-database defined name to the new extent -----
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Calibration!`$A`$2:`$N`$8"
    }
}

# --- Reselect the inserted rows, matching the saved selection -------------
$ws.Activate() | Out-Null
$ws.Range("A4:N8").Select() | Out-Null
